# Generate Report for Handback
#
# Renames the previously-handed-back file
#   a66e5554-e7d1-42ca-acc7-7540ff62b6d0.md  ->  190e3830-54a1-484c-9d22-5345a0c66512.md
# (updating its timestamps/xliff hashes) and appends a second, newly
# handed-back file
#   785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md
# as a new row on every sheet (Overview, zh-cn, de-de), keeping the
# tables/hyperlinks/dimensions in sync.

$wb = $excel.ActiveWorkbook

function Set-Text($ws, $addr, $text) {
    # Leading apostrophe forces literal-text entry so values that look like
    # booleans/numbers/dates ("True", "False", "", "2016-09-03 09:04:34", ...)
    # stay shared-string text cells instead of being coerced to t="b" etc.
    $ws.Range($addr).Value = "'" + $text
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

Set-Text $ws "A2" "190e3830-54a1-484c-9d22-5345a0c66512.md"
Set-Text $ws "B2" "e2e\190e3830-54a1-484c-9d22-5345a0c66512.md"
Set-Text $ws "C2" ".md"
Set-Text $ws "E2" "Handed back: in sync with en-US"
Set-Text $ws "F2" "Handed back: in sync with en-US"
Set-Text $ws "G2" "2016-09-03 09:04:34"

Set-Text $ws "A3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md"
Set-Text $ws "B3" "e2e\785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md"
Set-Text $ws "C3" ".md"
Set-Text $ws "E3" "Handed back: in sync with en-US"
Set-Text $ws "F3" "Handed back: in sync with en-US"
Set-Text $ws "G3" "2016-09-03 09:04:34"

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b7c31c12e259bfd6dae7382c2e6677fd9a1d6d8/e2e/190e3830-54a1-484c-9d22-5345a0c66512.md", $null, $null, "e2e\190e3830-54a1-484c-9d22-5345a0c66512.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b7c31c12e259bfd6dae7382c2e6677fd9a1d6d8/e2e/785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md", $null, $null, "e2e\785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

Set-Text $ws "A2" "190e3830-54a1-484c-9d22-5345a0c66512.md"
Set-Text $ws "B2" ".md"
Set-Text $ws "C2" "Handed back: in sync with en-US"
Set-Text $ws "D2" "e2e"
Set-Text $ws "E2" "ht"
Set-Text $ws "F2" "False"
Set-Text $ws "G2" "190e3830-54a1-484c-9d22-5345a0c66512.b54c89817eaced4dc7354d1fb6eb9bc98b3114ff.zh-cn.xlf"
Set-Text $ws "H2" "2016-09-03 09:04:29"
Set-Text $ws "I2" "190e3830-54a1-484c-9d22-5345a0c66512.md"
Set-Text $ws "J2" "190e3830-54a1-484c-9d22-5345a0c66512.b54c89817eaced4dc7354d1fb6eb9bc98b3114ff.zh-cn.xlf"
Set-Text $ws "K2" "2016-09-03 09:04:46"
Set-Text $ws "L2" ""
Set-Text $ws "M2" "True"
Set-Text $ws "N2" ""
Set-Text $ws "O2" "False"
Set-Text $ws "P2" ""

Set-Text $ws "A3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md"
Set-Text $ws "B3" ".md"
Set-Text $ws "C3" "Handed back: in sync with en-US"
Set-Text $ws "D3" "e2e"
Set-Text $ws "E3" "ht"
Set-Text $ws "F3" "True"
Set-Text $ws "G3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.46d8de1e06f1c2183f25e19e88d17b0552b85b97.zh-cn.xlf"
Set-Text $ws "H3" "2016-09-03 09:04:29"
Set-Text $ws "I3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md"
Set-Text $ws "J3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.46d8de1e06f1c2183f25e19e88d17b0552b85b97.zh-cn.xlf"
Set-Text $ws "K3" "2016-09-03 09:04:46"
Set-Text $ws "L3" ""
Set-Text $ws "M3" "True"
Set-Text $ws "N3" ""
Set-Text $ws "O3" "False"
Set-Text $ws "P3" ""

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a993fe8a430bdd23ed532bfaf39d7fe15ee6320/e2e/190e3830-54a1-484c-9d22-5345a0c66512.md", $null, $null, "190e3830-54a1-484c-9d22-5345a0c66512.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a993fe8a430bdd23ed532bfaf39d7fe15ee6320/e2e/190e3830-54a1-484c-9d22-5345a0c66512.md", $null, $null, "190e3830-54a1-484c-9d22-5345a0c66512.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a993fe8a430bdd23ed532bfaf39d7fe15ee6320/e2e/785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md", $null, $null, "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a993fe8a430bdd23ed532bfaf39d7fe15ee6320/e2e/785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md", $null, $null, "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

Set-Text $ws "A2" "190e3830-54a1-484c-9d22-5345a0c66512.md"
Set-Text $ws "B2" ".md"
Set-Text $ws "C2" "Handed back: in sync with en-US"
Set-Text $ws "D2" "e2e"
Set-Text $ws "E2" "ht"
Set-Text $ws "F2" "False"
Set-Text $ws "G2" "190e3830-54a1-484c-9d22-5345a0c66512.b54c89817eaced4dc7354d1fb6eb9bc98b3114ff.de-de.xlf"
Set-Text $ws "H2" "2016-09-03 09:04:34"
Set-Text $ws "I2" "190e3830-54a1-484c-9d22-5345a0c66512.md"
Set-Text $ws "J2" "190e3830-54a1-484c-9d22-5345a0c66512.b54c89817eaced4dc7354d1fb6eb9bc98b3114ff.de-de.xlf"
Set-Text $ws "K2" "2016-09-03 09:04:53"
Set-Text $ws "L2" ""
Set-Text $ws "M2" "True"
Set-Text $ws "N2" ""
Set-Text $ws "O2" "False"
Set-Text $ws "P2" ""

Set-Text $ws "A3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md"
Set-Text $ws "B3" ".md"
Set-Text $ws "C3" "Handed back: in sync with en-US"
Set-Text $ws "D3" "e2e"
Set-Text $ws "E3" "ht"
Set-Text $ws "F3" "True"
Set-Text $ws "G3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.46d8de1e06f1c2183f25e19e88d17b0552b85b97.de-de.xlf"
Set-Text $ws "H3" "2016-09-03 09:04:34"
Set-Text $ws "I3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md"
Set-Text $ws "J3" "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.46d8de1e06f1c2183f25e19e88d17b0552b85b97.de-de.xlf"
Set-Text $ws "K3" "2016-09-03 09:04:53"
Set-Text $ws "L3" ""
Set-Text $ws "M3" "True"
Set-Text $ws "N3" ""
Set-Text $ws "O3" "False"
Set-Text $ws "P3" ""

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/33fa0d1819a315b0614e2924dbbf7dc44efeb297/e2e/190e3830-54a1-484c-9d22-5345a0c66512.md", $null, $null, "190e3830-54a1-484c-9d22-5345a0c66512.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/33fa0d1819a315b0614e2924dbbf7dc44efeb297/e2e/190e3830-54a1-484c-9d22-5345a0c66512.md", $null, $null, "190e3830-54a1-484c-9d22-5345a0c66512.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/33fa0d1819a315b0614e2924dbbf7dc44efeb297/e2e/785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md", $null, $null, "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/33fa0d1819a315b0614e2924dbbf7dc44efeb297/e2e/785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md", $null, $null, "785a7ff2-c3c5-4c0a-a63e-d82900cce3d3.md") | Out-Null

Write-Output "Generate Report for Handback: done"
